$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.437.30"
$ws.Range("D3").Value = "1.864.90"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.55"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4777"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3761"
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07330"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9350"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  +4.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07822"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.864.32"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.556"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.57"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008885"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "27.504.72"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.939"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.34"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.46"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.966"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08894"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.330"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7538"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.597"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.733"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.120"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02033"
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05265"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.986"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5316"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.079"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.573"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.65"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4800"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.656"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.88"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.27"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06077"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9187"
$ws.Range("E51").Value = "  +3.48%  "
